$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for the corresponding rows, per the
# repulled data / mean calculation.
$values = @{
    2  = -2
    4  = -2
    5  = 1
    6  = -1
    8  = -2
    9  = -1
    10 = -4
    11 = 1
    12 = -3
    13 = -5
    14 = -1
    15 = -8
    16 = 1
    17 = -1
    18 = -3
    19 = 3
    20 = 7
    21 = 2
    22 = -1
    23 = -2
    24 = 4
    25 = -2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
